$d = $word.ActiveDocument

# --- Edit 1: empty paragraph gets justification removed and a new run "f" added ---
$p1 = $d.Paragraphs(6)
$p1.Format.Alignment = 0

# Source the exact run formatting (rStyle/rFonts/i/iCs/color+theme/sz/szCs) from the
# run immediately preceding this paragraph (the final "." of the previous paragraph),
# which already carries the desired rPr, then overwrite its text with "f".
$srcEnd = $d.Paragraphs(5).Range.End
$srcRun = $d.Range($srcEnd - 2, $srcEnd - 1)

$insPoint = $p1.Range.Start
$collapsed = $d.Range($insPoint, $insPoint)
$collapsed.FormattedText = $srcRun.FormattedText
$newChar = $d.Range($insPoint, $insPoint + 1)
$newChar.Text = "f"

# --- Edit 2 & 3: colorize "o problema do porquê não deu para criar." in red (EE0000) ---
# Paragraph index 37 ("Date: ...") -- keeps the trailing
# " Se estiver em branco, também tem que ir para lá" uncolored.
$p37 = $d.Paragraphs(37)

$a = $p37.Range.Duplicate
$null = $a.Find.Execute("o problema do ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$a.Font.Color = 238

$b = $p37.Range.Duplicate
$b.Start = $a.End
$null = $b.Find.Execute("porquê", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$b.Font.Color = 238

$c = $p37.Range.Duplicate
$c.Start = $b.End
$null = $c.Find.Execute(" não deu para criar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Font.Color = 238

$e = $p37.Range.Duplicate
$e.Start = $c.End
$null = $e.Find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$e.Font.Color = 238

# Paragraph index 39 ("Vencimento: ...") -- trailing text is
# ". Se estiver em branco, colocar o último dia útil do ano." and stays uncolored.
$p39 = $d.Paragraphs(39)

$g = $p39.Range.Duplicate
$null = $g.Find.Execute("o problema do ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$g.Font.Color = 238

$h = $p39.Range.Duplicate
$h.Start = $g.End
$null = $h.Find.Execute("porquê", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$h.Font.Color = 238

$i = $p39.Range.Duplicate
$i.Start = $h.End
$null = $i.Find.Execute(" não deu para criar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$i.Font.Color = 238
